$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ22949828",
    "summ23185337",
    "summ23432741",
    "summ23666181",
    "summ23903626",
    "summ24149658",
    "summ24382934",
    "summ24619721",
    "summ24850926"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
